# ---------------------------------------------------------------------------
# Updates the crypto price table (rows 2-51) to the refreshed coinranking.com
# snapshot: new Price/Volume(1h) figures for every row, and a one-row shift of
# the KuCoinToken..LEO block (rows 7-17) caused by GateToken re-entering the
# ranking above KuCoinToken.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name (B) / Link (C) column updates - plain text, no numeric coercion risk.
$textUpdates = [ordered]@{
    "B7" = "GateToken"
    "C7" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "B8" = "KuCoinToken"
    "C8" = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "B9" = "MXToken"
    "C9" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "B10" = "LiechtensteinCryptoassetsExchange"
    "C10" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B13" = "BitrueCoin"
    "C13" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "B16" = "TigerCash"
    "C16" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Price (D) / Volume(1h) (E) column updates. These look numeric/percent, so force
# the cell format to Text first - otherwise Excel auto-coerces "302.89" to a number
# and "-1.45%" to a percentage, losing the literal text representation (and, for
# values like "0.001660", the significant trailing zero).
$numericLookingUpdates = [ordered]@{
    "D2" = "302.89"
    "E2" = "-1.45%"
    "D3" = "35.65"
    "E3" = "-1.77%"
    "D4" = "5.038"
    "E4" = "-1.27%"
    "D5" = "0.07897"
    "E5" = "-2.95%"
    "D6" = "1.851"
    "E6" = "-5.48%"
    "D7" = "4.106"
    "E7" = "-2.19%"
    "D8" = "7.786"
    "E8" = "0.19%"
    "D9" = "0.9191"
    "E9" = "-1.04%"
    "D10" = "0.1343"
    "E10" = "-0.75%"
    "D11" = "0.1901"
    "E11" = "-1.40%"
    "D12" = "0.09094"
    "E12" = "-1.81%"
    "D13" = "0.03475"
    "E13" = "-3.58%"
    "D14" = "0.09835"
    "E14" = "-0.20%"
    "D15" = "0.001408"
    "E15" = "-0.51%"
    "D16" = "0.006142"
    "E16" = "6.22%"
    "D17" = "3.714"
    "E17" = "4.45%"
    "E18" = "11.97%"
    "E19" = "0.02%"
    "E20" = "3.03%"
    "D21" = "5.165"
    "E21" = "5.74%"
    "D22" = "0.2192"
    "E22" = "-8.98%"
    "D23" = "0.04409"
    "E23" = "-2.32%"
    "E24" = "1.73%"
    "D25" = "0.004617"
    "E25" = "-5.42%"
    "D26" = "0.0001301"
    "E26" = "4.93%"
    "D27" = "0.0004445"
    "E27" = "0.09%"
    "E39" = "-3.72%"
    "D40" = "0.05082"
    "E40" = "2.72%"
    "D41" = "0.007618"
    "E41" = "-0.22%"
    "D42" = "0.01016"
    "E42" = "-8.47%"
    "D43" = "0.1344"
    "E43" = "-2.64%"
    "D44" = "0.002152"
    "E44" = "2.47%"
    "D45" = "0.01018"
    "E45" = "-4.01%"
    "D46" = "0.00006183"
    "E46" = "-4.21%"
    "E47" = "0.08%"
    "D48" = "65.22"
    "E48" = "0.85%"
    "D49" = "0.001660"
    "E49" = "39.45%"
    "D50" = "0.00002102"
    "E50" = "0.08%"
    "D51" = "0.0002002"
    "E51" = "0.08%"
}
foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
}

